# live_trading_results.xlsx update
# Trade #50 (MarketMaking) closes via early_exit; Trade #83 (MarketMaking) opens.
# Updates the Summary / Strategy Status roll-ups and appends/edits the trade logs
# on "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a literal text value into a cell without Excel's automatic
# date/time/number inference turning it into something else (e.g. turning
# "2026-02-17" into a date serial). We temporarily force the cell to Text
# format, assign the literal value, then clear the format again so the
# resulting cell carries no explicit style - matching a freshly authored row.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.42   # Current Capital
$summary.Range("B4").Value = 0.21      # Total P&L $
$summary.Range("B6").Value = 50        # Total Trades
$summary.Range("B7").Value = 23        # Winning Trades
$summary.Range("B9").Value = 46        # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.42     # Capital
$status.Range("D5").Value = 17         # Trades
$status.Range("E5").Value = 0.1        # P&L $
$status.Range("F5").Value = 0.42       # P&L %
$status.Range("G5").Value = 52.94      # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet
#   Row 51 = Trade #50, closing out (early_exit)
#   Row 84 = new Trade #83, opening
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G51").Value = 0.84
$allTrades.Range("H51").Value = "CLOSED"
$allTrades.Range("I51").Value = 2.439
$allTrades.Range("J51").Value = 0.02
$allTrades.Range("K51").Value = 100.42
Set-TextValue $allTrades.Range("L51") "early_exit"
$allTrades.Range("M51").Value = 0.14

Set-TextValue $allTrades.Range("B84") "2026-02-17"
Set-TextValue $allTrades.Range("C84") "20:53:32"
$allTrades.Range("A84").Value = 83
Set-TextValue $allTrades.Range("D84") "MarketMaking"
Set-TextValue $allTrades.Range("E84") "DOWN"
$allTrades.Range("F84").Value = 0.82
$allTrades.Range("H84").Value = "OPEN"
$allTrades.Range("I84").Value = 0
$allTrades.Range("J84").Value = 0
$allTrades.Range("K84").Value = 100.3984370824165
$allTrades.Range("M84").Value = 0
$allTrades.Range("N84").Value = 0
$allTrades.Range("O84").Value = 0
$allTrades.Range("P84").Value = 0.6
Set-TextValue $allTrades.Range("Q84") "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet
#   Row 18 = Trade #50, closing out (early_exit)
#   Row 51 = new Trade #83, opening
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Range("G18").Value = 0.84
$mm.Range("H18").Value = "CLOSED"
$mm.Range("I18").Value = 2.439
$mm.Range("J18").Value = 0.02
$mm.Range("K18").Value = 100.42
Set-TextValue $mm.Range("P18") "early_exit"
$mm.Range("Q18").Value = 0.14

Set-TextValue $mm.Range("B51") "2026-02-17"
Set-TextValue $mm.Range("C51") "20:53:32"
$mm.Range("A51").Value = 83
Set-TextValue $mm.Range("D51") "MarketMaking"
Set-TextValue $mm.Range("E51") "DOWN"
$mm.Range("F51").Value = 0.82
$mm.Range("H51").Value = "OPEN"
$mm.Range("I51").Value = 0
$mm.Range("J51").Value = 0
$mm.Range("K51").Value = 100.3984370824165
$mm.Range("L51").Value = 0
$mm.Range("M51").Value = 0
$mm.Range("N51").Value = 0.6
Set-TextValue $mm.Range("O51") "Normal spread capture: 19600 bps"
$mm.Range("Q51").Value = 0
